$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($text -eq "Step 8:") {
        $p.Range.Delete()
        break
    }
}
